# Adds a new "2022-Q3" quarterly sheet ahead of "2022-Q2", fills it with the
# new fund-holding data, inserts a matching summary row on "总计", and fixes
# a stale header label on "2021-Q2".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right before "2022-Q2" so the tab
#    order becomes: 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3,
#    2021-Q2, 2021-Q1, 2020-Q4
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Match the bold / bordered / centered look used by the header row and the
# leading index column on every other quarterly sheet.
$headerRng = $q3.Range("B1:H1")
$headerRng.Font.Bold = $true
$headerRng.Borders.LineStyle = 1
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160

$indexRng = $q3.Range("A2:A3")
$indexRng.Font.Bold = $true
$indexRng.Borders.LineStyle = 1
$indexRng.HorizontalAlignment = -4108
$indexRng.VerticalAlignment = -4160

# Header row (row 1) - same layout as every other quarterly sheet.
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# The B:G columns hold text (exactly like the other quarterly sheets), so
# force text formatting before assigning the values.
$q3.Range("B2:G3").NumberFormat = "@"

# Row 2 - 920002 中金精选股票A
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "920002"
$q3.Range("C2").Value = "中金精选股票A"
$q3.Range("D2").Value = "2.95"
$q3.Range("E2").Value = "82.28"
$q3.Range("F2").Value = "2.38"
$q3.Range("G2").Value = "0.0702"
$q3.Range("H2").Value = 10

# Row 3 - 920922 中金精选股票C
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "920922"
$q3.Range("C3").Value = "中金精选股票C"
$q3.Range("D3").Value = "0.11"
$q3.Range("E3").Value = "82.28"
$q3.Range("F3").Value = "2.38"
$q3.Range("G3").Value = "0.0026"
$q3.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2) "总计" (summary) sheet - insert a new row for 2022-Q3 above the
#    existing 2022-Q2 row, pushing everything else down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Row 3 (the old "2022-Q2" row, now shifted down) carries the correct
# formatting for column A (bold + border) - copy it onto the new A2 so
# the inserted row matches the sheet's look.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.07

# ---------------------------------------------------------------------
# 3) "2021-Q2" sheet - the "基金规模" column header used to read
#    "基金金额"; fix the stale label.
# ---------------------------------------------------------------------
$q2_2021 = $wb.Worksheets.Item("2021-Q2")
$q2_2021.Range("D1").Value = "基金规模"
